# Apply updated ifoCAST matched-error values (shifted one horizon, with newly
# evaluated vintages appended) to rows 2-24, columns B:K.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$row2 = New-Object "object[,]" 1,10
$row2[0,0] = 9.634600668846847
$row2[0,1] = -8.622337597701534
$row2[0,2] = -0.8522339323542271
$row2[0,3] = 0.478975912873543
$row2[0,4] = -1.892429927382574
$row2[0,5] = -1.228164494743756
$row2[0,6] = -0.4586921403196634
$row2[0,7] = -0.5513017252472258
$row2[0,8] = 0.1696252062987764
$row2[0,9] = -0.5393267183150783
$ws.Range("B2:K2").Value2 = $row2

# Row 3
$row3 = New-Object "object[,]" 1,10
$row3[0,0] = -10.87214043731293
$row3[0,1] = -3.102036771965619
$row3[0,2] = -1.770826926737849
$row3[0,3] = -4.142232766993966
$row3[0,4] = -3.477967334355148
$row3[0,5] = -2.708494979931056
$row3[0,6] = -2.801104564858618
$row3[0,7] = -2.080177633312616
$row3[0,8] = -2.78912955792647
$row3[0,9] = -2.297023565603993
$ws.Range("B3:K3").Value2 = $row3

# Row 4
$row4 = New-Object "object[,]" 1,10
$row4[0,0] = -10.57118710282291
$row4[0,1] = -9.239977257595143
$row4[0,2] = -11.61138309785126
$row4[0,3] = -10.94711766521244
$row4[0,4] = -10.17764531078835
$row4[0,5] = -10.27025489571591
$row4[0,6] = -9.549327964169908
$row4[0,7] = -10.25827988878376
$row4[0,8] = -9.766173896461286
$row4[0,9] = -9.2963691491883
$ws.Range("B4:K4").Value2 = $row4

# Row 5
$row5 = New-Object "object[,]" 1,10
$row5[0,0] = 8.931284542542326
$row5[0,1] = 6.55987870228621
$row5[0,2] = 7.224144134925027
$row5[0,3] = 7.99361648934912
$row5[0,4] = 7.901006904421558
$row5[0,5] = 8.62193383596756
$row5[0,6] = 7.912981911353705
$row5[0,7] = 8.405087903676183
$row5[0,8] = 8.87489265094917
$row5[0,9] = 9.187717845914616
$ws.Range("B5:K5").Value2 = $row5

# Row 6
$row6 = New-Object "object[,]" 1,10
$row6[0,0] = -1.579792844692939
$row6[0,1] = -0.9155274120541215
$row6[0,2] = -0.1460550576300293
$row6[0,3] = -0.2386646425575917
$row6[0,4] = 0.4822622889884105
$row6[0,5] = -0.2266896356254442
$row6[0,6] = 0.2654163566970329
$row6[0,7] = 0.73522110397002
$row6[0,8] = 1.048046298935466
$row6[0,9] = 0.169643717683886
$ws.Range("B6:K6").Value2 = $row6

# Row 7
$row7 = New-Object "object[,]" 1,10
$row7[0,0] = -3.063097883853514
$row7[0,1] = -2.293625529429422
$row7[0,2] = -2.386235114356984
$row7[0,3] = -1.665308182810982
$row7[0,4] = -2.374260107424837
$row7[0,5] = -1.88215411510236
$row7[0,6] = -1.412349367829372
$row7[0,7] = -1.099524172863926
$row7[0,8] = -1.977926754115506
$row7[0,9] = -2.707838769293975
$ws.Range("B7:K7").Value2 = $row7

# Row 8
$row8 = New-Object "object[,]" 1,10
$row8[0,0] = 1.146404457093908
$row8[0,1] = 1.053794872166346
$row8[0,2] = 1.774721803712348
$row8[0,3] = 1.065769879098493
$row8[0,4] = 1.55787587142097
$row8[0,5] = 2.027680618693958
$row8[0,6] = 2.340505813659403
$row8[0,7] = 1.462103232407824
$row8[0,8] = 0.7321912172293545
$row8[0,9] = 2.394370313618982
$ws.Range("B8:K8").Value2 = $row8

# Row 9
$row9 = New-Object "object[,]" 1,10
$row9[0,0] = 1.114969050580547
$row9[0,1] = 1.835895982126549
$row9[0,2] = 1.126944057512694
$row9[0,3] = 1.619050049835171
$row9[0,4] = 2.088854797108159
$row9[0,5] = 2.401679992073604
$row9[0,6] = 1.523277410822025
$row9[0,7] = 0.7933653956435556
$row9[0,8] = 2.455544492033183
$row9[0,9] = 1.836325203637939
$ws.Range("B9:K9").Value2 = $row9

# Row 10
$row10 = New-Object "object[,]" 1,10
$row10[0,0] = -0.2055599550297054
$row10[0,1] = -0.9145118796435601
$row10[0,2] = -0.422405887321083
$row10[0,3] = 0.04739885995190407
$row10[0,4] = 0.36022405491735
$row10[0,5] = -0.5181785263342299
$row10[0,6] = -1.248090541512699
$row10[0,7] = 0.4140885548769285
$row10[0,8] = -0.2051307335183153
$row10[0,9] = 0.03217303010139827
$ws.Range("B10:K10").Value2 = $row10

# Row 11
$row11 = New-Object "object[,]" 1,10
$row11[0,0] = -0.3319173426099191
$row11[0,1] = 0.160188649712558
$row11[0,2] = 0.6299933969855451
$row11[0,3] = 0.942818591950991
$row11[0,4] = 0.06441601069941108
$row11[0,5] = -0.6654960044790579
$row11[0,6] = 0.9966830919105695
$row11[0,7] = 0.3774638035153257
$row11[0,8] = 0.6147675671350392
$row11[0,9] = 0.5114901972596275
$ws.Range("B11:K11").Value2 = $row11

# Row 12
$row12 = New-Object "object[,]" 1,10
$row12[0,0] = 0.0645136089032002
$row12[0,1] = 0.5343183561761873
$row12[0,2] = 0.8471435511416332
$row12[0,3] = -0.03125903010994671
$row12[0,4] = -0.7611710452884157
$row12[0,5] = 0.9010080511012117
$row12[0,6] = 0.2817887627059679
$row12[0,7] = 0.5190925263256815
$row12[0,8] = 0.4158151564502698
$row12[0,9] = -0.3757832708791649
$ws.Range("B12:K12").Value2 = $row12

# Row 13
$row13 = New-Object "object[,]" 1,10
$row13[0,0] = 0.7947373931749101
$row13[0,1] = 1.107562588140356
$row13[0,2] = 0.2291600068887761
$row13[0,3] = -0.5007520082896928
$row13[0,4] = 1.161427088099934
$row13[0,5] = 0.5422077997046907
$row13[0,6] = 0.7795115633244043
$row13[0,7] = 0.6762341934489926
$row13[0,8] = -0.1153642338804421
$row13[0,9] = 0.4951246737870189
$ws.Range("B13:K13").Value2 = $row13

# Row 14
$row14 = New-Object "object[,]" 1,10
$row14[0,0] = 0.2721075919818648
$row14[0,1] = -0.6062949892697151
$row14[0,2] = -1.336207004448184
$row14[0,3] = 0.3259720919414433
$row14[0,4] = -0.2932471964538005
$row14[0,5] = -0.05594343283408693
$row14[0,6] = -0.1592208027094986
$row14[0,7] = -0.9508192300389333
$row14[0,8] = -0.3403303223714723
$row14[0,9] = -0.5018523531907899
$ws.Range("B14:K14").Value2 = $row14

# Row 15
$row15 = New-Object "object[,]" 1,9
$row15[0,0] = -0.5062156293670936
$row15[0,1] = -1.236127644545562
$row15[0,2] = 0.4260514518440648
$row15[0,3] = -0.193167836551179
$row15[0,4] = 0.04413592706853459
$row15[0,5] = -0.05914144280687711
$row15[0,6] = -0.8507398701363118
$row15[0,7] = -0.2402509624688508
$row15[0,8] = -0.4017729932881683
$ws.Range("B15:J15").Value2 = $row15
$ws.Range("K15:K15").ClearContents()

# Row 16
$row16 = New-Object "object[,]" 1,8
$row16[0,0] = -0.8823411384658664
$row16[0,1] = 0.779837957923761
$row16[0,2] = 0.1606186695285172
$row16[0,3] = 0.3979224331482308
$row16[0,4] = 0.2946450632728191
$row16[0,5] = -0.4969533640566156
$row16[0,6] = 0.1135355436108454
$row16[0,7] = -0.04798648720847212
$ws.Range("B16:I16").Value2 = $row16
$ws.Range("J16:K16").ClearContents()

# Row 17
$row17 = New-Object "object[,]" 1,7
$row17[0,0] = 0.6316605674913157
$row17[0,1] = 0.0124412790960719
$row17[0,2] = 0.2497450427157855
$row17[0,3] = 0.1464676728403738
$row17[0,4] = -0.6451307544890609
$row17[0,5] = -0.03464184682159993
$row17[0,6] = -0.1961638776409175
$ws.Range("B17:H17").Value2 = $row17
$ws.Range("I17:K17").ClearContents()

# Row 18
$row18 = New-Object "object[,]" 1,6
$row18[0,0] = -0.1450047099080831
$row18[0,1] = 0.0922990537116305
$row18[0,2] = -0.0109783161637812
$row18[0,3] = -0.8025767434932158
$row18[0,4] = -0.1920878358257549
$row18[0,5] = -0.3536098666450724
$ws.Range("B18:G18").Value2 = $row18
$ws.Range("H18:K18").ClearContents()

# Row 19
$row19 = New-Object "object[,]" 1,5
$row19[0,0] = 0.5929585102377013
$row19[0,1] = 0.4896811403622896
$row19[0,2] = -0.3019172869671451
$row19[0,3] = 0.3085716207003159
$row19[0,4] = 0.1470495898809984
$ws.Range("B19:F19").Value2 = $row19
$ws.Range("G19:K19").ClearContents()

# Row 20
$row20 = New-Object "object[,]" 1,4
$row20[0,0] = 0.2093232598268204
$row20[0,1] = -0.5822751675026142
$row20[0,2] = 0.02821374016484672
$row20[0,3] = -0.1333082906544708
$ws.Range("B20:E20").Value2 = $row20
$ws.Range("F20:K20").ClearContents()

# Row 21
$row21 = New-Object "object[,]" 1,3
$row21[0,0] = -0.4103003096576026
$row21[0,1] = 0.2001885980098584
$row21[0,2] = 0.03866656719054083
$ws.Range("B21:D21").Value2 = $row21
$ws.Range("E21:K21").ClearContents()

# Row 22
$row22 = New-Object "object[,]" 1,2
$row22[0,0] = -0.1056739417364731
$row22[0,1] = -0.2671959725557906
$ws.Range("B22:C22").Value2 = $row22
$ws.Range("D22:K22").ClearContents()

# Row 23
$row23 = New-Object "object[,]" 1,1
$row23[0,0] = 0.3451339801314955
$ws.Range("B23:B23").Value2 = $row23
$ws.Range("C23:K23").ClearContents()

# Row 24
$ws.Range("B24:K24").ClearContents()

